$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.043.03"
$ws.Range("E2").Value = "  +0.38%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.384.55"
$ws.Range("E3").Value = "  -1.01%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.25"
$ws.Range("E5").Value = "  +1.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.85"
$ws.Range("E6").Value = "  -2.43%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -0.66%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.106"
$ws.Range("E9").Value = "  +0.64%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.67"
$ws.Range("E10").Value = "  -0.50%  "

# Row 11
$ws.Range("E11").Value = "  +1.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.343"
$ws.Range("E12").Value = "  -3.14%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.47"
$ws.Range("E13").Value = "  -3.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.797.66"
$ws.Range("E14").Value = "  -1.29%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.927.97"
$ws.Range("E15").Value = "  +0.35%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000138"
$ws.Range("E16").Value = "  +0.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.381.96"
$ws.Range("E17").Value = "  -1.45%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.11"
$ws.Range("E18").Value = "  -1.80%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.48"
$ws.Range("E19").Value = "  +1.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.47"
$ws.Range("E20").Value = "  -2.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.13"
$ws.Range("E23").Value = "  -3.61%  "

# Row 24
$ws.Range("E24").Value = "  +0.25%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.45"
$ws.Range("E26").Value = "  -2.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.38"
$ws.Range("E27").Value = "  +0.64%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.81"
$ws.Range("E28").Value = "  +2.23%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0761"
$ws.Range("E29").Value = "  -1.36%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.95"
$ws.Range("E30").Value = "  +0.98%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.08"
$ws.Range("E31").Value = "  +0.95%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("E32").Value = "  +11.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.398"
$ws.Range("E33").Value = "  -2.10%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.19"
$ws.Range("E34").Value = "  -2.38%  "

# Row 35
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.33"
$ws.Range("E36").Value = "  +1.93%  "

# Row 37
$ws.Range("E37").Value = "  +0.10%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.15"
$ws.Range("E38").Value = "  -0.98%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.59"
$ws.Range("E39").Value = "  -0.69%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "319.73"
$ws.Range("E40").Value = "  +0.04%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.61"
$ws.Range("E41").Value = "  -1.85%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "145.22"
$ws.Range("E42").Value = "  +4.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.54"
$ws.Range("E43").Value = "  -3.71%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0970"
$ws.Range("E44").Value = "  +0.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.73"
$ws.Range("E45").Value = "  +0.49%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0511"
$ws.Range("E46").Value = "  -0.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.571"
$ws.Range("E47").Value = "  -1.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0218"
$ws.Range("E48").Value = "  -2.27%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.07"
$ws.Range("E49").Value = "  +0.31%  "

# Row 50
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.54"
$ws.Range("E50").Value = "  -1.50%  "

# Row 51
$ws.Range("B51").Value = "ZEEBU"
$ws.Range("C51").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.68"
$ws.Range("E51").Value = "  +0.14%  "
